$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 85.521736
$ws.Range("I11").Value = 85.521736
$ws.Range("K11").Value = 85.521736
$ws.Range("M11").Value = 54.478264
$ws.Range("H19").Value = 1023.625
$ws.Range("I19").Value = 575
$ws.Range("J19").Value = 1173.1666
$ws.Range("K19").Value = 575
$ws.Range("L19").Value = 1173.1666
$ws.Range("M19").Value = -400
$ws.Range("N19").Value = -1523.1666
$ws.Range("H61").Value = 1005.25
$ws.Range("I61").Value = 673.6667
$ws.Range("K61").Value = 2021.0001
$ws.Range("M61").Value = -1849.0001
$ws.Range("H118").Value = 880.2
$ws.Range("I118").Value = 1010.3333
$ws.Range("K118").Value = 3030.9999
$ws.Range("M118").Value = -1373.9999
$ws.Range("H137").Value = 1923.3334
$ws.Range("I137").Value = 1737.0769
$ws.Range("J137").Value = 2615.1428
$ws.Range("K137").Value = 5211.2307
$ws.Range("L137").Value = 7845.428400000001
$ws.Range("M137").Value = -2661.2307
$ws.Range("N137").Value = -12945.4284
$ws.Range("H141").Value = 3140.7693
$ws.Range("I141").Value = 3251
$ws.Range("K141").Value = 9753
$ws.Range("M141").Value = -4573

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2628.6924
$ws.Range("I74").Value = 3084.1
$ws.Range("K74").Value = 3084.1
$ws.Range("M74").Value = -2210.1
$ws.Range("H77").Value = 2628.6924
$ws.Range("I77").Value = 3084.1
$ws.Range("K77").Value = 15420.5
$ws.Range("M77").Value = -11052.5
$ws.Range("H102").Value = 3676.4736
$ws.Range("I102").Value = 1682.3077
$ws.Range("J102").Value = 7997.1665
$ws.Range("K102").Value = 1682.3077
$ws.Range("L102").Value = 7997.1665
$ws.Range("M102").Value = -60.30770000000007
$ws.Range("N102").Value = -11241.1665
$ws.Range("H110").Value = 645.5
$ws.Range("I110").Value = 645.5
$ws.Range("K110").Value = 645.5
$ws.Range("M110").Value = 1399.5
$ws.Range("H122").Value = 1808.9
$ws.Range("I122").Value = 1455.1482
$ws.Range("K122").Value = 4365.444600000001
$ws.Range("M122").Value = -1915.444600000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 11484.777
$ws.Range("I22").Value = 14654
$ws.Range("J22").Value = 392.5
$ws.Range("K22").Value = 14654
$ws.Range("L22").Value = 392.5
$ws.Range("M22").Value = -14481
$ws.Range("N22").Value = -738.5
$ws.Range("H94").Value = 1093.9
$ws.Range("I94").Value = 1142.4839
$ws.Range("J94").Value = 926.55554
$ws.Range("K94").Value = 1142.4839
$ws.Range("L94").Value = 926.55554
$ws.Range("M94").Value = -691.4838999999999
$ws.Range("N94").Value = -1828.55554
$ws.Range("H107").Value = 3891.125
$ws.Range("I107").Value = 2632.9092
$ws.Range("J107").Value = 6659.2
$ws.Range("K107").Value = 2632.9092
$ws.Range("L107").Value = 6659.2
$ws.Range("M107").Value = -712.9092000000001
$ws.Range("N107").Value = -10499.2

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 4273.88
$ws.Range("I7").Value = 281.64706
$ws.Range("K7").Value = 281.64706
$ws.Range("M7").Value = -168.64706
$ws.Range("H16").Value = 1971.2858
$ws.Range("I16").Value = 1799.8334
$ws.Range("K16").Value = 1799.8334
$ws.Range("M16").Value = -1512.8334
$ws.Range("H29").Value = 11500
$ws.Range("I29").Value = 9750
$ws.Range("K29").Value = 9750
$ws.Range("M29").Value = -9457
$ws.Range("H31").Value = 4882.684
$ws.Range("I31").Value = 2774.6667
$ws.Range("J31").Value = 5277.9375
$ws.Range("K31").Value = 2774.6667
$ws.Range("L31").Value = 5277.9375
$ws.Range("M31").Value = -2479.6667
$ws.Range("N31").Value = -5867.9375
$ws.Range("H34").Value = 4882.684
$ws.Range("I34").Value = 2774.6667
$ws.Range("J34").Value = 5277.9375
$ws.Range("K34").Value = 2774.6667
$ws.Range("L34").Value = 5277.9375
$ws.Range("M34").Value = -2572.6667
$ws.Range("N34").Value = -5681.9375
$ws.Range("H62").Value = 6499
$ws.Range("J62").Value = 6499
$ws.Range("L62").Value = 6499
$ws.Range("N62").Value = -7747
$ws.Range("H65").Value = 6499
$ws.Range("J65").Value = 6499
$ws.Range("L65").Value = 32495
$ws.Range("N65").Value = -38735
$ws.Range("H80").Value = 64999.5
$ws.Range("J80").Value = 64999.5
$ws.Range("L80").Value = 64999.5
$ws.Range("N80").Value = -67245.5
$ws.Range("H83").Value = 64999.5
$ws.Range("J83").Value = 64999.5
$ws.Range("L83").Value = 194998.5
$ws.Range("N83").Value = -206230.5
$ws.Range("H94").Value = 1465.5454
$ws.Range("I94").Value = 267.66666
$ws.Range("K94").Value = 267.66666
$ws.Range("M94").Value = 183.33334
$ws.Range("H113").Value = 1971.2858
$ws.Range("I113").Value = 1799.8334
$ws.Range("K113").Value = 1799.8334
$ws.Range("M113").Value = 370.1666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 75
$ws.Range("I21").Value = 75
$ws.Range("K21").Value = 225
$ws.Range("M21").Value = -52
$ws.Range("H108").Value = 6249.9
$ws.Range("I108").Value = 306.5
$ws.Range("J108").Value = 30023.5
$ws.Range("K108").Value = 919.5
$ws.Range("L108").Value = 90070.5
$ws.Range("M108").Value = 1960.5
$ws.Range("N108").Value = -95830.5
$ws.Range("H110").Value = 99.5
$ws.Range("I110").Value = 99.5
$ws.Range("K110").Value = 298.5
$ws.Range("M110").Value = 3791.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 56000.2
$ws.Range("J26").Value = 56000.2
$ws.Range("L26").Value = 56000.2
$ws.Range("N26").Value = -56560.2
$ws.Range("H50").Value = 56000.2
$ws.Range("J50").Value = 56000.2
$ws.Range("L50").Value = 56000.2
$ws.Range("N50").Value = -56996.2
$ws.Range("H52").Value = 43329.332
$ws.Range("J52").Value = 43329.332
$ws.Range("L52").Value = 43329.332
$ws.Range("N52").Value = -43847.332
$ws.Range("H102").Value = 2058.3333
$ws.Range("I102").Value = 1518.1818
$ws.Range("K102").Value = 1518.1818
$ws.Range("M102").Value = 103.8181999999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2008.4706
$ws.Range("I46").Value = 871.4286
$ws.Range("J46").Value = 2804.4
$ws.Range("K46").Value = 871.4286
$ws.Range("L46").Value = 2804.4
$ws.Range("M46").Value = -683.4286
$ws.Range("N46").Value = -3180.4
$ws.Range("H55").Value = 876.48486
$ws.Range("J55").Value = 1826.9286
$ws.Range("L55").Value = 1826.9286
$ws.Range("N55").Value = -2172.9286
$ws.Range("H61").Value = 4234.75
$ws.Range("I61").Value = 3665.4333
$ws.Range("K61").Value = 3665.4333
$ws.Range("M61").Value = -3463.4333
$ws.Range("H68").Value = 4303.1035
$ws.Range("I68").Value = 4005.7058
$ws.Range("J68").Value = 4724.4165
$ws.Range("K68").Value = 4005.7058
$ws.Range("L68").Value = 4724.4165
$ws.Range("M68").Value = -3256.7058
$ws.Range("N68").Value = -6222.4165
$ws.Range("H69").Value = 600039800
$ws.Range("J69").Value = 600039800
$ws.Range("L69").Value = 600039800
$ws.Range("N69").Value = -600041422
$ws.Range("H71").Value = 4303.1035
$ws.Range("I71").Value = 4005.7058
$ws.Range("J71").Value = 4724.4165
$ws.Range("K71").Value = 20028.529
$ws.Range("L71").Value = 23622.0825
$ws.Range("M71").Value = -16284.529
$ws.Range("N71").Value = -31110.0825
$ws.Range("H72").Value = 600039800
$ws.Range("J72").Value = 600039800
$ws.Range("L72").Value = 1800119400
$ws.Range("N72").Value = -1800127512
$ws.Range("H80").Value = 68450
$ws.Range("J80").Value = 68450
$ws.Range("L80").Value = 68450
$ws.Range("N80").Value = -70696
$ws.Range("H82").Value = 1333.6842
$ws.Range("I82").Value = 999.2727
$ws.Range("J82").Value = 1793.5
$ws.Range("K82").Value = 999.2727
$ws.Range("L82").Value = 1793.5
$ws.Range("M82").Value = -638.2727
$ws.Range("N82").Value = -2515.5
$ws.Range("H83").Value = 68450
$ws.Range("J83").Value = 68450
$ws.Range("L83").Value = 205350
$ws.Range("N83").Value = -216582
$ws.Range("H85").Value = 1333.6842
$ws.Range("I85").Value = 999.2727
$ws.Range("J85").Value = 1793.5
$ws.Range("K85").Value = 999.2727
$ws.Range("L85").Value = 1793.5
$ws.Range("M85").Value = 248.7273
$ws.Range("N85").Value = -4289.5
$ws.Range("H96").Value = 42999.5
$ws.Range("J96").Value = 42999.5
$ws.Range("L96").Value = 42999.5
$ws.Range("N96").Value = -48491.5
$ws.Range("H113").Value = 4234.75
$ws.Range("I113").Value = 3665.4333
$ws.Range("K113").Value = 3665.4333
$ws.Range("M113").Value = -1495.4333
$ws.Range("H132").Value = 4559
$ws.Range("J132").Value = 4559
$ws.Range("L132").Value = 13677
$ws.Range("N132").Value = -18737

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 23285.715
$ws.Range("J49").Value = 23285.715
$ws.Range("L49").Value = 23285.715
$ws.Range("N49").Value = -23745.715
$ws.Range("H54").Value = 13651.059
$ws.Range("I54").Value = 982.1429000000001
$ws.Range("J54").Value = 22519.3
$ws.Range("K54").Value = 982.1429000000001
$ws.Range("L54").Value = 22519.3
$ws.Range("M54").Value = -462.1429000000001
$ws.Range("N54").Value = -23559.3
$ws.Range("H126").Value = 2683.2144
$ws.Range("I126").Value = 2107.1428
$ws.Range("J126").Value = 3259.2856
$ws.Range("K126").Value = 6321.428400000001
$ws.Range("L126").Value = 9777.856800000001
$ws.Range("M126").Value = -3851.428400000001
$ws.Range("N126").Value = -14717.8568
$ws.Range("H132").Value = 4266.579
$ws.Range("I132").Value = 5112.1665
$ws.Range("K132").Value = 15336.4995
$ws.Range("M132").Value = -12806.4995
$ws.Range("H136").Value = 1462.8223
$ws.Range("I136").Value = 946.8919
$ws.Range("K136").Value = 2840.6757
$ws.Range("M136").Value = -290.6756999999998
